# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a set of rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of worksheet row number -> new DAMSLTag / DialogAct values
$updates = @{
    6   = @("b",  "Acknowledge (Backchannel)")
    9   = @("sd", "Statement-non-opinion")
    34  = @("sv", "Statement-opinion")
    39  = @("aa", "Agree/Accept")
    47  = @("aa", "Agree/Accept")
    52  = @("b",  "Acknowledge (Backchannel)")
    80  = @("aa", "Agree/Accept")
    87  = @("sv", "Statement-opinion")
    90  = @("ba", "Appreciation")
    91  = @("b",  "Acknowledge (Backchannel)")
    101 = @("sd", "Statement-non-opinion")
    106 = @("sv", "Statement-opinion")
    108 = @("%",  "Uninterpretable")
    112 = @("ba", "Appreciation")
    116 = @("%",  "Uninterpretable")
    117 = @("ba", "Appreciation")
    119 = @("sd", "Statement-non-opinion")
    132 = @("aa", "Agree/Accept")
    133 = @("sd", "Statement-non-opinion")
    139 = @("sd", "Statement-non-opinion")
    141 = @("b",  "Acknowledge (Backchannel)")
    148 = @("%",  "Uninterpretable")
    154 = @("sv", "Statement-opinion")
    155 = @("sv", "Statement-opinion")
    161 = @("b",  "Acknowledge (Backchannel)")
    166 = @("ba", "Appreciation")
    172 = @("sd", "Statement-non-opinion")
    183 = @("ba", "Appreciation")
    184 = @("sd", "Statement-non-opinion")
    198 = @("sd", "Statement-non-opinion")
    201 = @("sv", "Statement-opinion")
    204 = @("ba", "Appreciation")
    228 = @("ba", "Appreciation")
    229 = @("sv", "Statement-opinion")
    234 = @("b",  "Acknowledge (Backchannel)")
    239 = @("sd", "Statement-non-opinion")
    241 = @("aa", "Agree/Accept")
    243 = @("ba", "Appreciation")
    251 = @("sv", "Statement-opinion")
    256 = @("sd", "Statement-non-opinion")
    262 = @("b",  "Acknowledge (Backchannel)")
    263 = @("b",  "Acknowledge (Backchannel)")
    266 = @("aa", "Agree/Accept")
    270 = @("sd", "Statement-non-opinion")
    271 = @("aa", "Agree/Accept")
    287 = @("%",  "Uninterpretable")
    288 = @("aa", "Agree/Accept")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]   # Column I - DAMSLTag
    $ws.Cells.Item($row, 10).Value = $values[1]  # Column J - DialogAct
}
